$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pelada (Saturday match) results typed in below the existing table:
# player name (shared string) + Vitorias, Empate, Derrotas, Gols, Partidas,
# Tarde de Vitoria, La barca, Craque do Dia, Gols Sofridos, Melhor Goleiro.
$newRows = @(@(134,'Athos',3,2,4,3,1,0,0,0,0,0),@(135,'Marcos',3,2,4,1,1,0,0,0,0,0),@(136,'Cabeleira',3,2,4,2,1,0,0,0,0,0),@(137,'Corinthiano',3,2,4,1,1,0,0,0,0,0),@(138,'Cristiano',3,2,4,0,1,0,0,0,0,0),@(139,'Leandrinho',3,2,3,2,1,0,0,0,0,0),@(140,'Digão',3,2,3,2,1,0,0,0,0,0),@(141,'Jorge',3,2,3,2,1,0,0,0,0,0),@(142,'Fabinho',3,2,3,1,1,0,0,0,0,0),@(143,'Marcelão',3,2,3,1,1,0,0,0,0,0),@(144,'Guinha',4,2,3,1,1,1,0,0,0,0),@(145,'Miqueias',4,2,3,1,1,1,0,0,0,0),@(146,'Senna',4,2,3,1,1,1,0,0,0,0),@(147,'Vander',4,2,3,1,1,1,0,0,0,0),@(148,'Leah',4,2,3,1,1,1,0,1,0,0),@(149,'Romario',2,2,2,0,0,0,1,0,0,0),@(150,'Juscielio',2,2,2,1,0,0,1,0,0,0),@(151,'Du',2,2,2,1,0,0,1,0,0,0),@(152,'Peixe',2,2,2,2,0,0,1,0,0,0),@(153,'Coxinha',2,2,2,0,0,0,1,0,0,0),@(154,'Lucian',0,2,0,0,1,0,0,0,3,0),@(155,'Alan',6,3,5,0,1,1,0,0,11,1),@(156,'Matheus',5,3,6,0,1,0,1,0,9,0))

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $newRows[$i]
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    for ($j = 0; $j -lt 10; $j++) {
        $ws.Cells.Item($rowNum, 3 + $j).Value = $r[2 + $j]
    }
}

# Leave the selection where the user's cursor ended up after typing the
# last new row (the header-row freeze is untouched).
$ws.Range("A157").Select()
